# Commit: Updated symbol list on Wed Jan 18 10:07:41 UTC 2023 with GitHub Actions
#
# The scraper re-ran and refreshed three columns of the crypto table for every
# data row (rows 2-51): D = Price, E = Volume(1h), G = Hora (hour snapshot).
# Column A (index), B (Coin), C (Link) and F (Data/date) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: @(sheet-row, sheet-column, new literal text value).
# Columns: D=4 (Price), E=5 (Volume 1h), G=7 (Hora).
# Values are written as text (NumberFormat "@") so they stay literal strings -
# e.g. "0.07670" keeps its trailing zero and "10" stays text, not a number -
# matching the original inlineStr cells in the workbook.
$updates = @(
    @(2, 4, "300.57"),
    @(2, 5, "0.00%"),
    @(2, 7, "10"),
    @(3, 4, "32.67"),
    @(3, 5, "3.65%"),
    @(3, 7, "10"),
    @(4, 4, "4.951"),
    @(4, 5, "-2.89%"),
    @(4, 7, "10"),
    @(5, 4, "0.07670"),
    @(5, 5, "-2.14%"),
    @(5, 7, "10"),
    @(6, 4, "1.918"),
    @(6, 5, "-17.70%"),
    @(6, 7, "10"),
    @(7, 4, "7.826"),
    @(7, 5, "0.22%"),
    @(7, 7, "10"),
    @(8, 4, "3.803"),
    @(8, 5, "-0.98%"),
    @(8, 7, "10"),
    @(9, 4, "0.9192"),
    @(9, 5, "0.09%"),
    @(9, 7, "10"),
    @(10, 4, "0.1752"),
    @(10, 5, "-0.58%"),
    @(10, 7, "10"),
    @(11, 4, "0.07761"),
    @(11, 5, "2.65%"),
    @(11, 7, "10"),
    @(12, 4, "0.08526"),
    @(12, 5, "-6.92%"),
    @(12, 7, "10"),
    @(13, 4, "0.03168"),
    @(13, 5, "3.00%"),
    @(13, 7, "10"),
    @(14, 5, "-0.07%"),
    @(14, 7, "10"),
    @(15, 4, "0.001511"),
    @(15, 5, "-0.26%"),
    @(15, 7, "10"),
    @(16, 4, "0.005884"),
    @(16, 5, "0.87%"),
    @(16, 7, "10"),
    @(17, 5, "-0.23%"),
    @(17, 7, "10"),
    @(18, 4, "2.154"),
    @(18, 5, "-4.11%"),
    @(18, 7, "10"),
    @(19, 4, "0.3352"),
    @(19, 5, "2.45%"),
    @(19, 7, "10"),
    @(20, 4, "0.1326"),
    @(20, 5, "-0.77%"),
    @(20, 7, "10"),
    @(21, 4, "4.272"),
    @(21, 5, "6.67%"),
    @(21, 7, "10"),
    @(22, 4, "0.1992"),
    @(22, 5, "11.37%"),
    @(22, 7, "10"),
    @(23, 4, "0.04519"),
    @(23, 5, "-1.69%"),
    @(23, 7, "10"),
    @(24, 4, "0.001221"),
    @(24, 5, "-2.17%"),
    @(24, 7, "10"),
    @(25, 4, "0.004404"),
    @(25, 7, "10"),
    @(26, 4, "0.0001252"),
    @(26, 5, "0.10%"),
    @(26, 7, "10"),
    @(27, 7, "10"),
    @(28, 7, "10"),
    @(29, 7, "10"),
    @(30, 7, "10"),
    @(31, 7, "10"),
    @(32, 7, "10"),
    @(33, 7, "10"),
    @(34, 7, "10"),
    @(35, 7, "10"),
    @(36, 7, "10"),
    @(37, 7, "10"),
    @(38, 7, "10"),
    @(39, 4, "0.01700"),
    @(39, 5, "-3.98%"),
    @(39, 7, "10"),
    @(40, 4, "0.04694"),
    @(40, 5, "-2.81%"),
    @(40, 7, "10"),
    @(41, 4, "0.007496"),
    @(41, 5, "4.44%"),
    @(41, 7, "10"),
    @(42, 4, "0.1350"),
    @(42, 7, "10"),
    @(43, 4, "0.002333"),
    @(43, 5, "6.50%"),
    @(43, 7, "10"),
    @(44, 4, "0.01045"),
    @(44, 5, "1.77%"),
    @(44, 7, "10"),
    @(45, 4, "0.00006238"),
    @(45, 5, "-1.64%"),
    @(45, 7, "10"),
    @(46, 5, "0.07%"),
    @(46, 7, "10"),
    @(47, 5, "10.44%"),
    @(47, 7, "10"),
    @(48, 7, "10"),
    @(49, 5, "0.07%"),
    @(49, 7, "10"),
    @(50, 5, "0.07%"),
    @(50, 7, "10"),
    @(51, 7, "10"),
)

foreach ($u in $updates) {
    $r = $u[0]
    $c = $u[1]
    $v = $u[2]
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $v
}
